# Apply updates to the Day Sale report:
#  - Update the "current balance" (الرصيد الحالي) values for a few items
#    that now have matching/returned stock recorded ("0:0" -> "1:0").
#  - Refresh the generation timestamp shown at the bottom of the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H11").Value = "1:0"
$ws.Range("H16").Value = "1:0"
$ws.Range("H18").Value = "1:0"
$ws.Range("H23").Value = "1:0"

$ws.Range("A28").Value = "Tuesday, 19 August, 2025 10:49 AM"
